# Contest 53 CSK vs PBKS (row 62) & Contest 54 KKR vs RR (row 63)
# Fill in the raw score inputs (columns E, H, K, N, Q, T) for both
# contests. The dependent VLOOKUP/RANK formulas in D, G, J, P (and, for
# row 63, M and S as well) recalc automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 62 - Contest 53: CSK vs PBKS
$ws.Range("E62").Value = 80
$ws.Range("H62").Value = 100
$ws.Range("K62").Value = 60
$ws.Range("N62").Value = 40
$ws.Range("Q62").Value = 0
$ws.Range("T62").Value = 40

# Row 63 - Contest 54: KKR vs RR
$ws.Range("E63").Value = 80
$ws.Range("H63").Value = 60
$ws.Range("K63").Value = 100
$ws.Range("N63").Value = 40
$ws.Range("Q63").Value = 0
$ws.Range("T63").Value = 20

# M62 and S62 tied in RANK (N62 and T62 are both 40), so the author
# manually overwrote those two formula cells with the averaged
# tie-break value instead of leaving the RANK/VLOOKUP formula in place.
$ws.Range("M62").Value = -22.5
$ws.Range("S62").Value = -22.5
